# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
# These two sheets mirror the same data, so the same set of updates is applied
# to each.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value (applies to column F on both sheets).
$updates = @{
    "F3"  = 1782
    "F5"  = 397
    "F6"  = 234
    "F8"  = 194
    "F11" = 25
    "F12" = 95
    "F13" = 248
    "F18" = 73
    "F19" = 263
    "F20" = 36
    "F21" = 447
    "F22" = 348
    "F24" = 54
    "F25" = 28
    "F26" = 35
    "F27" = 786
    "F28" = 2588
    "F31" = 512
    "F32" = 842
    "F34" = 445
    "F35" = 250
    "F36" = 375
    "F37" = 443
    "F38" = 579
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
